$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$old1 = $cellA1.Value2
$new1 = $old1.Replace("1000 Bs = 9.84 = 41456.0 pesos", "1000 Bs = 9.81 = 41219.44 pesos")
$new1 = $new1.Replace("41456.0 pesos = 9.8 = 977.25 Bs", "41219.44 pesos = 9.74 = 971.68 Bs")
$cellA1.Value = $new1

# --- Sheet "tasas": update the rate cells N10, O10, N12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 101.981
$wsTasas.Range("O10").Value = 4203.6
$wsTasas.Range("N12").Value = 4229.99
